$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map numeric month values (1-12) to Spanish month abbreviations, matching
# the "Mes" column text used throughout the report.
$monthNames = @{
    1  = "Ene."
    2  = "Feb."
    3  = "Mar."
    4  = "Abr."
    5  = "May."
    6  = "Jun."
    7  = "Jul."
    8  = "Ago."
    9  = "Sep."
    10 = "Oct."
    11 = "Nov."
    12 = "Dic."
}

for ($row = 6; $row -le 80; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $monthNum = [int]$cell.Value2
    $cell.Value = $monthNames[$monthNum]
}
